$d = $word.ActiveDocument

# --- 1. Insert the "Meta description" paragraph right after the H1 title ---
$titlePara = $d.Paragraphs(1)
$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Dolphin Gold with Stellar Jackpots. Play for free and enjoy exciting gameplay, stunning graphics, and triple jackpots.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Range.InsertXML($metaXml)

# --- 2. Remove the trailing duplicate "Play ... | Review" bold paragraph ---
# (the title paragraph at index 1 must stay; only delete the later duplicate)
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 2; $i--) {
    $para = $d.Paragraphs($i)
    $t = $para.Range.Text.Trim()
    if ($t -eq "Play Dolphin Gold with Stellar Jackpots Free | Review") {
        $para.Range.Delete()
        break
    }
}

# --- 3. Update the trailing italic meta-description-like paragraph to the new Prompt text ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$r = $lastPara.Range
$contentRange = $d.Range($r.Start, $r.End - 1)
$contentRange.Text = "Prompt: Create a cartoon-style feature image for Dolphin Gold with Stellar Jackpots that features a happy Maya warrior wearing glasses. The image should have an underwater theme with the dolphin and gold elements incorporated into the background. The Maya warrior should be holding a treasure chest and smiling at the viewer. Use bright colors and bold lines to make the image pop and attract attention to the game's exciting features. The image should convey the idea of adventure and treasure while also showcasing the game's playful and enjoyable aspects."

Write-Output "Done"
